$d = $word.ActiveDocument

# Locate the paragraph that holds the literal "${ KategoriProyek }" placeholder
# (there is a second, differently-structured "${KategoriProyek}" occurrence
# inside the table further down, which must stay untouched).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`$`{ KategoriProyek }`r") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph '`${ KategoriProyek }'"
}

$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="46ECB3D6" w14:textId="2EA547BE" w:rsidR="006525C6" w:rsidRDefault="00956FE1" w:rsidP="003D7C7D"><w:pPr><w:spacing w:before="2" w:after="0" w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00956FE1"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>KategoriProyek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace the whole paragraph (range includes the paragraph mark) so the
# stray w:proofErr[@type="gramStart"] that sits just before the run sequence
# is dropped along with it.
$target.Range.InsertXML($newParaXml)

Write-Output "Done. New text: $($d.Content.Text)"
